# Auto-generated edit script applying scheduled-runner price updates
# to the Asura_Profits workbook (sheet-by-sheet: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 26400
$ws.Range("J57").Value = 26400
$ws.Range("L57").Value = 79200
$ws.Range("N57").Value = -80198
# Row 57 done

$ws.Range("H64").Value = 3355.1035
$ws.Range("I64").Value = 3379.8
$ws.Range("K64").Value = 3379.8
$ws.Range("M64").Value = -3131.8
# Row 64 done

$ws.Range("H67").Value = 3355.1035
$ws.Range("I67").Value = 3379.8
$ws.Range("K67").Value = 3379.8
$ws.Range("M67").Value = -2521.8
# Row 67 done

$ws.Range("H88").Value = 3140.4
$ws.Range("I88").Value = 5000
$ws.Range("J88").Value = 2933.7778
$ws.Range("K88").Value = 5000
$ws.Range("L88").Value = 2933.7778
$ws.Range("M88").Value = -4594
$ws.Range("N88").Value = -3745.7778
# Row 88 done

$ws.Range("H91").Value = 3140.4
$ws.Range("I91").Value = 5000
$ws.Range("J91").Value = 2933.7778
$ws.Range("K91").Value = 5000
$ws.Range("L91").Value = 2933.7778
$ws.Range("M91").Value = -3596
$ws.Range("N91").Value = -5741.7778
# Row 91 done

$ws.Range("H116").Value = 8698130
$ws.Range("J116").Value = 2656.1428
$ws.Range("L116").Value = 2656.1428
$ws.Range("N116").Value = -9540.1428
# Row 116 done

$ws.Range("H129").Value = 1155.762
$ws.Range("J129").Value = 1342.1212
$ws.Range("L129").Value = 4026.3636
$ws.Range("N129").Value = -14026.3636
# Row 129 done

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 20833.334
$ws.Range("J34").Value = 20833.334
$ws.Range("L34").Value = 20833.334
$ws.Range("N34").Value = -21375.334
# Row 34 done

$ws.Range("H45").Value = 1073.65
$ws.Range("I45").Value = 1042.6875
$ws.Range("J45").Value = 1197.5
$ws.Range("K45").Value = 1042.6875
$ws.Range("L45").Value = 1197.5
$ws.Range("M45").Value = -665.6875
$ws.Range("N45").Value = -1951.5
# Row 45 done

$ws.Range("H74").Value = 898.1842
$ws.Range("I74").Value = 735
$ws.Range("K74").Value = 735
$ws.Range("M74").Value = 139
# Row 74 done

$ws.Range("H77").Value = 898.1842
$ws.Range("I77").Value = 735
$ws.Range("K77").Value = 3675
$ws.Range("M77").Value = 693
# Row 77 done

$ws.Range("H88").Value = 2962
$ws.Range("I88").Value = 2481.2
$ws.Range("J88").Value = 3262.5
$ws.Range("K88").Value = 2481.2
$ws.Range("L88").Value = 3262.5
$ws.Range("M88").Value = -2075.2
$ws.Range("N88").Value = -4074.5
# Row 88 done

$ws.Range("H91").Value = 2962
$ws.Range("I91").Value = 2481.2
$ws.Range("J91").Value = 3262.5
$ws.Range("K91").Value = 2481.2
$ws.Range("L91").Value = 3262.5
$ws.Range("M91").Value = -1077.2
$ws.Range("N91").Value = -6070.5
# Row 91 done

$ws.Range("H132").Value = 3572.082
$ws.Range("I132").Value = 3856.5227
$ws.Range("K132").Value = 11569.5681
$ws.Range("M132").Value = -9039.5681
# Row 132 done

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H69").Value = 32000
$ws.Range("J69").Value = 32000
$ws.Range("L69").Value = 32000
$ws.Range("N69").Value = -33622
# Row 69 done

$ws.Range("H72").Value = 32000
$ws.Range("J72").Value = 32000
$ws.Range("L72").Value = 96000
$ws.Range("N72").Value = -104112
# Row 72 done

$ws.Range("H134").Value = 2437.6191
$ws.Range("I134").Value = 2128.0322
$ws.Range("J134").Value = 3310.0908
$ws.Range("K134").Value = 6384.096600000001
$ws.Range("L134").Value = 9930.2724
$ws.Range("M134").Value = -3849.096600000001
$ws.Range("N134").Value = -15000.2724
# Row 134 done

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 118.652176
$ws.Range("I7").Value = 84
$ws.Range("J7").Value = 150.41667
$ws.Range("K7").Value = 84
$ws.Range("L7").Value = 150.41667
$ws.Range("M7").Value = 29
$ws.Range("N7").Value = -376.41667
# Row 7 done

$ws.Range("H22").Value = 409.8095
$ws.Range("I22").Value = 235.25
$ws.Range("J22").Value = 968.4
$ws.Range("K22").Value = 235.25
$ws.Range("L22").Value = 968.4
$ws.Range("M22").Value = 114.75
$ws.Range("N22").Value = -1668.4
# Row 22 done

$ws.Range("H35").Value = 1392.4
$ws.Range("I35").Value = 1392.4
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1392.4
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -1098.4
$ws.Range("N35").ClearContents()
# Row 35 done

$ws.Range("H132").Value = 437646.53
$ws.Range("I132").Value = 521220.38
$ws.Range("K132").Value = 1563661.14
$ws.Range("M132").Value = -1561131.14
# Row 132 done

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 540.3333
$ws.Range("I114").Value = 389.44446
$ws.Range("J114").Value = 766.6667
$ws.Range("K114").Value = 1168.33338
$ws.Range("L114").Value = 2300.0001
$ws.Range("M114").Value = 2085.66662
$ws.Range("N114").Value = -8808.000100000001
# Row 114 done

$ws.Range("H131").Value = 903.49493
$ws.Range("J131").Value = 942.34784
$ws.Range("L131").Value = 2827.04352
$ws.Range("N131").Value = -12907.04352
# Row 131 done

$ws.Range("H134").Value = 3640.2927
$ws.Range("I134").Value = 1926.7142
$ws.Range("J134").Value = 5439.55
$ws.Range("K134").Value = 5780.142599999999
$ws.Range("L134").Value = 16318.65
$ws.Range("M134").Value = -710.1425999999992
$ws.Range("N134").Value = -26458.65
# Row 134 done

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 50002
$ws.Range("J20").Value = 50003
$ws.Range("L20").Value = 50003
$ws.Range("N20").Value = -50493
# Row 20 done

$ws.Range("H102").Value = 2056.0857
$ws.Range("I102").Value = 1554.16
$ws.Range("J102").Value = 3310.9
$ws.Range("K102").Value = 1554.16
$ws.Range("L102").Value = 3310.9
$ws.Range("M102").Value = 67.83999999999992
$ws.Range("N102").Value = -6554.9
# Row 102 done

$ws.Range("H123").Value = 1808392.5
$ws.Range("J123").Value = 1808392.5
$ws.Range("L123").Value = 1808392.5
$ws.Range("N123").Value = -1813292.5
# Row 123 done

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 859.1579
$ws.Range("I22").Value = 977.7143
$ws.Range("J22").Value = 790
$ws.Range("K22").Value = 977.7143
$ws.Range("L22").Value = 790
$ws.Range("M22").Value = -682.7143
$ws.Range("N22").Value = -1380
# Row 22 done

$ws.Range("H27").Value = 859.1579
$ws.Range("I27").Value = 977.7143
$ws.Range("J27").Value = 790
$ws.Range("K27").Value = 977.7143
$ws.Range("L27").Value = 790
$ws.Range("M27").Value = -870.7143
$ws.Range("N27").Value = -1004
# Row 27 done

$ws.Range("H32").Value = 2144.2856
$ws.Range("I32").Value = 2144.2856
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2144.2856
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1827.2856
$ws.Range("N32").ClearContents()
# Row 32 done

$ws.Range("H61").Value = 35067.332
$ws.Range("I61").Value = 41080.8
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 41080.8
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -40878.8
$ws.Range("N61").Value = -5404
# Row 61 done

$ws.Range("H113").Value = 35067.332
$ws.Range("I113").Value = 41080.8
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 41080.8
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -38910.8
$ws.Range("N113").Value = -9340
# Row 113 done

$ws.Range("H132").Value = 3758.3447
$ws.Range("I132").Value = 3630.2354
$ws.Range("J132").Value = 3939.8333
$ws.Range("K132").Value = 10890.7062
$ws.Range("L132").Value = 11819.4999
$ws.Range("M132").Value = -8360.706200000001
$ws.Range("N132").Value = -16879.4999
# Row 132 done

$ws.Range("H136").Value = 25899500
$ws.Range("I136").Value = 33334590
$ws.Range("K136").Value = 100003770
$ws.Range("M136").Value = -100001220
# Row 136 done

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 8288.416999999999
$ws.Range("I126").Value = 10175.667
$ws.Range("J126").Value = 2626.6667
$ws.Range("K126").Value = 30527.001
$ws.Range("L126").Value = 7880.000100000001
$ws.Range("M126").Value = -28057.001
$ws.Range("N126").Value = -12820.0001
# Row 126 done

$ws.Range("H132").Value = 1574.5818
$ws.Range("I132").Value = 986.72974
$ws.Range("K132").Value = 2960.18922
$ws.Range("M132").Value = -430.1892200000002
# Row 132 done

$ws.Range("H136").Value = 1391.8
$ws.Range("I136").Value = 1132.6072
$ws.Range("K136").Value = 3397.8216
$ws.Range("M136").Value = -847.8215999999998
# Row 136 done
